# Atualização do código do produto do Módulo 1
# "Casos acumulados" table: switch thousands-separators from commas to
# periods (e.g. "10,132,849" -> "10.132.849") for every country row.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "10,132,849";  New = "10.132.849"  },
    @{ Old = "11,861,161";  New = "11.861.161"  },
    @{ Old = "37,511,921";  New = "37.511.921"  },
    @{ Old = "99,361,338";  New = "99.361.338"  },
    @{ Old = "38,997,490";  New = "38.997.490"  },
    @{ Old = "38,437,756";  New = "38.437.756"  },
    @{ Old = "45,040,074";  New = "45.040.074"  },
    @{ Old = "26,727,644";  New = "26.727.644"  },
    @{ Old = "33,803,572";  New = "33.803.572"  },
    @{ Old = "34,571,873";  New = "34.571.873"  },
    @{ Old = "24,225,459";  New = "24.225.459"  },
    @{ Old = "13,980,340";  New = "13.980.340"  },
    @{ Old = "17,004,714";  New = "17.004.714"  },
    @{ Old = "24,940,968";  New = "24.940.968"  },
    @{ Old = "103,436,829"; New = "103.436.829" },
    @{ Old = "11,624,000";  New = "11.624.000"  }
)

# Walk every paragraph (the numbers each live alone in a table-cell
# paragraph) and, when its text matches one of the old values, set the
# paragraph range's Text directly. Assigning Range.Text (rather than
# Find/Replace) keeps the run's existing formatting/xml:space intact and
# only swaps the literal digits, exactly mirroring the authored edit.
foreach ($para in $d.Paragraphs) {
    $text = $para.Range.Text
    foreach ($r in $replacements) {
        if ($text -like "*$($r.Old)*") {
            $para.Range.Text = $r.New
        }
    }
}
